# Updates the cryptos list (Sheet1) to reflect the latest scrape.
# Numeric-looking "Price" strings are written with a leading quote-prefix
# (the classic "force text" trick) so Excel keeps them as literal text
# (matching formats like "1.00" / "0.0467") instead of coercing them to
# numbers; the style is then reset to "Normal" so no stray quote-prefix
# cell format lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.320.73'
$ws.Range('E2').Value = '  +2.64%  '
$ws.Range('D3').Value = '3.234.15'
$ws.Range('E3').Value = '  +4.61%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '''576.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.73%  '
$ws.Range('D6').Value = '''154.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.14%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.226.64'
$ws.Range('E8').Value = '  +4.70%  '
$ws.Range('D9').Value = '''0.513'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.45%  '
$ws.Range('D10').Value = '''7.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.37%  '
$ws.Range('D11').Value = '''0.166'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.71%  '
$ws.Range('D12').Value = '''0.485'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.49%  '
$ws.Range('E13').Value = '  +5.81%  '
$ws.Range('E14').Value = '  +2.70%  '
$ws.Range('D15').Value = '3.764.04'
$ws.Range('E15').Value = '  +5.09%  '
$ws.Range('D16').Value = '66.430.66'
$ws.Range('E16').Value = '  +2.75%  '
$ws.Range('D17').Value = '''552.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +11.11%  '
$ws.Range('D18').Value = '3.236.28'
$ws.Range('E18').Value = '  +4.61%  '
$ws.Range('E19').Value = '  +3.13%  '
$ws.Range('D20').Value = '''7.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.69%  '
$ws.Range('D21').Value = '''14.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.24%  '
$ws.Range('E22').Value = '  +6.35%  '
$ws.Range('D23').Value = '''7.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.46%  '
$ws.Range('D24').Value = '''13.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.40%  '
$ws.Range('D25').Value = '''82.08'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.74%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '''9.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +16.03%  '
$ws.Range('D28').Value = '''2.92'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.57%  '
$ws.Range('D29').Value = '''2.27'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.45%  '
$ws.Range('D30').Value = '''27.82'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.67%  '
$ws.Range('D31').Value = '''2.76'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('D32').Value = '''0.999'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('E33').Value = '  +5.08%  '
$ws.Range('D34').Value = '''567.14'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.07%  '
$ws.Range('E35').Value = '  +3.56%  '
$ws.Range('D36').Value = '''6.45'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.32%  '
$ws.Range('D37').Value = '''0.0467'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +14.28%  '
$ws.Range('D38').Value = '''54.86'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.07%  '
$ws.Range('D39').Value = '''0.0875'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.10%  '
$ws.Range('D40').Value = '''3.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.43%  '
$ws.Range('E41').Value = '  +4.57%  '
$ws.Range('D42').Value = '3.127.50'
$ws.Range('E42').Value = '  +6.53%  '
$ws.Range('E43').Value = '  +2.41%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '''2.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.89%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '''0.274'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.20%  '
$ws.Range('D46').Value = '''27.18'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.31%  '
$ws.Range('D47').Value = '0.0₃0564'
$ws.Range('E47').Value = '  +3.31%  '
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('E49').Value = '  +4.07%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = '''2.26'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.36%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '''122.39'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.10%  '

Write-Output "Applied 97 cell updates"
